$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E16").Value = 338
$ws.Range("F16").Value = 150
$ws.Range("D17").Value = 120
$ws.Range("E17").Value = 330
$ws.Range("F17").Value = 107
$ws.Range("E18").Value = 667
$ws.Range("F18").Value = 189
$ws.Range("D19").Value = 402
$ws.Range("E19").Value = 979
$ws.Range("F19").Value = 343
$ws.Range("D20").Value = 175
$ws.Range("E20").Value = 433
$ws.Range("F20").Value = 127
$ws.Range("D21").Value = 104
$ws.Range("E21").Value = 231
$ws.Range("F21").Value = 63
$ws.Range("D22").Value = 137
$ws.Range("F22").Value = 53
$ws.Range("D23").Value = 148
$ws.Range("E23").Value = 380
$ws.Range("D24").Value = 137
$ws.Range("E24").Value = 338
$ws.Range("D25").Value = 315
$ws.Range("F25").Value = 150
$ws.Range("D26").Value = 465
$ws.Range("E26").Value = 1005
$ws.Range("F26").Value = 255
$ws.Range("E27").Value = 411
$ws.Range("F27").Value = 108
$ws.Range("D28").Value = 116
$ws.Range("E28").Value = 228
$ws.Range("D29").Value = 160
$ws.Range("E29").Value = 203
$ws.Range("D30").Value = 113
$ws.Range("E30").Value = 398
$ws.Range("E31").Value = 354
$ws.Range("F31").Value = 74
$ws.Range("D32").Value = 237
$ws.Range("E32").Value = 713
$ws.Range("F32").Value = 175
